$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before row 27 (pushes the existing rows 27..120
# down to 29..122, and the sheet's used range grows to A1:R122).
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(27).Insert()

# --- New row 27 -----------------------------------------------------------
$ws.Range("A27").Value = 9
$ws.Range("B27").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = "2022-12-19"
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = 100114002
$ws.Range("G27").Value = "Camote"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 790
$ws.Range("K27").Value = 17000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 17494
$ws.Range("N27").Value = "$/caja 18 kilos"
$ws.Range("O27").Value = "Perú"
$ws.Range("P27").Value = 972
$ws.Range("Q27").Value = 18
$ws.Range("R27").Value = "Hortaliza"

# --- New row 28 -----------------------------------------------------------
$ws.Range("A28").Value = 9
$ws.Range("B28").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = "2022-12-19"
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 100114002
$ws.Range("G28").Value = "Camote"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 520
$ws.Range("K28").Value = 18000
$ws.Range("L28").Value = 19000
$ws.Range("M28").Value = 18500
$ws.Range("N28").Value = "$/malla 18 kilos"
$ws.Range("O28").Value = "Perú"
$ws.Range("P28").Value = 1028
$ws.Range("Q28").Value = 18
$ws.Range("R28").Value = "Hortaliza"
